$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "type" column header (J2) ---
$ws.Range("J2").Value = "type (1:chq pt compte 2:have service for pt)"

# --- Fill in J column "type" values for existing rows 3-11 ---
$ws.Range("J3").Value = 2
$ws.Range("J4").Value = 2
$ws.Range("J5").Value = 1
$ws.Range("J6").Value = 1
$ws.Range("J7").Value = 1
$ws.Range("J8").Value = 1
$ws.Range("J9").Value = 1
$ws.Range("J10").Value = 1
$ws.Range("J11").Value = 1

# --- New match rows 12-35 ---
# Columns: row, B(player), C(player), F(score1), G(score2), H(game of), I(date), J(type)
$rows = @(
    @{r=12; b="W"; c="E"; f=11; g=7;  h=11; i="26/08/2024"; j=2},
    @{r=13; b="E"; c="W"; f=11; g=7;  h=11; i="26/08/2024"; j=2},
    @{r=14; b="E"; c="W"; f=12; g=10; h=11; i="26/08/2024"; j=2},
    @{r=15; b="E"; c="W"; f=11; g=3;  h=11; i="26/08/2024"; j=2},
    @{r=16; b="E"; c="W"; f=11; g=5;  h=11; i="26/08/2024"; j=2},
    @{r=17; b="E"; c="W"; f=11; g=6;  h=11; i="26/08/2024"; j=2},
    @{r=18; b="E"; c="W"; f=12; g=10; h=11; i="26/08/2024"; j=2},
    @{r=19; b="E"; c="W"; f=12; g=10; h=11; i="26/08/2024"; j=2},
    @{r=20; b="W"; c="E"; f=11; g=6;  h=11; i="26/08/2024"; j=2},
    @{r=21; b="E"; c="W"; f=15; g=13; h=11; i="28/08/2024"; j=1},
    @{r=22; b="E"; c="X"; f=11; g=6;  h=11; i="28/08/2024"; j=1},
    @{r=23; b="W"; c="X"; f=13; g=11; h=11; i="28/08/2024"; j=1},
    @{r=24; b="W"; c="E"; f=11; g=2;  h=11; i="28/08/2024"; j=1},
    @{r=25; b="E"; c="X"; f=11; g=6;  h=11; i="28/08/2024"; j=1},
    @{r=26; b="X"; c="W"; f=11; g=5;  h=11; i="28/08/2024"; j=1},
    @{r=27; b="W"; c="E"; f=11; g=9;  h=11; i="28/08/2024"; j=1},
    @{r=28; b="E"; c="X"; f=11; g=9;  h=11; i="28/08/2024"; j=1},
    @{r=29; b="X"; c="W"; f=11; g=9;  h=11; i="28/08/2024"; j=1},
    @{r=30; b="E"; c="W"; f=11; g=7;  h=11; i="28/08/2024"; j=1},
    @{r=31; b="X"; c="E"; f=11; g=5;  h=11; i="28/08/2024"; j=1},
    @{r=32; b="W"; c="X"; f=11; g=9;  h=11; i="28/08/2024"; j=1},
    @{r=33; b="W"; c="E"; f=11; g=6;  h=11; i="28/08/2024"; j=1},
    @{r=34; b="E"; c="X"; f=11; g=5;  h=11; i="28/08/2024"; j=1},
    @{r=35; b="W"; c="E"; f=12; g=10; h=11; i="28/08/2024"; j=1}
)

# "E" placeholder above stands for the accented player initial; map to the real text.
$playerMap = @{ "W" = "W"; "X" = "X"; "E" = [char]0x00C9 }

foreach ($row in $rows) {
    $r = $row.r
    $ws.Cells.Item($r, 2).Value = $playerMap[$row.b]
    $ws.Cells.Item($r, 3).Value = $playerMap[$row.c]

    $eFormula = "=IF(F" + $r + ">G" + $r + ",1,0)"
    $ws.Cells.Item($r, 5).Formula = $eFormula

    $ws.Cells.Item($r, 6).Value = $row.f
    $ws.Cells.Item($r, 7).Value = $row.g
    $ws.Cells.Item($r, 8).Value = $row.h
    $ws.Cells.Item($r, 9).Value = $row.i
    $ws.Cells.Item($r, 10).Value = $row.j

    $kFormula = '=IF(OR(OR(AND(OR(A' + $r + '=B' + $r + ',A' + $r + '=C' + $r + ',A' + $r + '=D' + $r + ',B' + $r + '=C' + $r + ',B' + $r + '=D' + $r + ',C' + $r + '=D' + $r + '),OR(A' + $r + '<>"",D' + $r + '<>"")),H' + $r + '>MAX(F' + $r + ':G' + $r + '),B' + $r + '=C' + $r + '),OR(AND(ISBLANK(A' + $r + ')=FALSE,ISNA(VLOOKUP(A' + $r + ',$M$3:$O$27,1,FALSE))),ISNA(VLOOKUP(B' + $r + ',$M$3:$O$27,1,FALSE)),ISNA(VLOOKUP(C' + $r + ',$M$3:$O$27,1,FALSE)),AND(ISBLANK(D' + $r + ')=FALSE,ISNA(VLOOKUP(D' + $r + ',$M$3:$O$27,1,FALSE))),OR(COUNTBLANK(A' + $r + ':D' + $r + ')=1,COUNTBLANK(A' + $r + ':D' + $r + ')=3))),"ERREUR","")'
    $ws.Cells.Item($r, 11).Formula = $kFormula
}

# --- View state: scrolled down a bit with a new selection (D33) ---
$ws.Range("D33").Select()
